# Rename the "SwateTemplateMetadata" sheet to "isa_template" and make it
# the active/selected tab (previously "Nanodrop_measurement" was active).

$wb = $excel.ActiveWorkbook

$metaSheet = $wb.Worksheets.Item("SwateTemplateMetadata")
$metaSheet.Name = "isa_template"

# Switch the active tab from the first sheet to the (renamed) second sheet.
$metaSheet.Activate()
